# Insert a new data row at row 40 (pushing rows 40..76 down to 41..77),
# and populate the new row 40 with the "Agrícola del Norte S.A. de Arica"
# Poroto verde record dated 44904.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 40, shifting existing data down.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new record's values.
$ws.Cells.Item(40, 1).Value  = 1
$ws.Cells.Item(40, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(40, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(40, 4).Value  = 44904
$ws.Cells.Item(40, 5).Value  = 15
$ws.Cells.Item(40, 6).Value  = 100112031
$ws.Cells.Item(40, 7).Value  = "Poroto verde"
$ws.Cells.Item(40, 8).Value  = "Sin especificar"
$ws.Cells.Item(40, 9).Value  = "Primera"
$ws.Cells.Item(40, 10).Value = 4000
$ws.Cells.Item(40, 11).Value = 500
$ws.Cells.Item(40, 12).Value = 600
$ws.Cells.Item(40, 13).Value = 550
$ws.Cells.Item(40, 14).Value = "`$/kilo"
$ws.Cells.Item(40, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(40, 16).Value = 550
$ws.Cells.Item(40, 17).Value = 1
$ws.Cells.Item(40, 18).Value = "Hortaliza"
